$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$objetivosText = @"
Apresentar aos alunos as bases teóricas e experimentais dos métodos instrumentais (quantitativos e qualitativos) de uso mais frequente na área química, incluindo o preparo de amostras e a criteriosa avaliação dos resultados analíticos. Ao final da disciplina, o aluno deve ser capaz de escolher e aplicar a metodologia mais adequada à solução dos problemas analíticos em geral, assim como interpretar resultados de análises químicas.
"@
$docentesValue = @"
2341641 - Maria da Rosa Capri
"@
$programaResumidoText = @"
Introdução à Análise Instrumental. Preparo de amostras. Métodos Espectroanalíticos: UV/Visível, Absorção Atômica, Emissão Atômica, Infravermelho. Métodos Eletroanalíticos: Potenciometria e Condutimetria. Métodos Cromatográficos: Cromatografia a Gás e Cromatografia Líquida de Alta Eficiência.
"@
$programaText = @"
1) Introdução à Análise Instrumental. Correlação entre métodos analíticos instrumentais e por via úmida. Preparo de amostras em meio sólido e em meios líquidos aquosos e não aquosos. Solubilização, digestão, fontes de energia aplicadas ao preparo, estabilização de amostras.
2) Introdução aos Métodos Espectroanalíticos: Natureza da energia radiante. Espectro eletromagnético. Interação da radiação com a matéria. Absorção seletiva. Absortividade. Lei de Beer-Lambert. Curvas analíticas. 
3) Introdução à Espectrofotometria no UV/Visível. Instrumentação. Aplicações e interpretação de resultados. Determinações simultâneas. Parte Experimental.
4) Introdução às Espectrometrias de Absorção e de Emissão Atômicas. Instrumentação. Interferências. Origem do espectro de emissão atômica. Fontes de atomização e de excitação. Calibração. Aplicações e interpretação de resultados. Parte Experimental.
5) Introdução à Espectroscopia no Infravermelho. Instrumentação. Interpretação de espectros. Aplicações. Parte Experimental.
6) Introdução aos Métodos Eletroanalíticos: Potenciometria e Condutimetria. Instrumentação.  Métodos diretos e indiretos. Aplicações e interpretação de resultados. Parte experimental.
7) Introdução aos Métodos Cromatográficos. Conceitos básicos dos métodos de separação. Fases móvel e estacionária. Cromatografia planar em papel e em camada delgada. Cromatografia em coluna: cromatografia a gás e cromatografia líquida de alta eficiência. Instrumentação. Aplicações e interpretação de resultados. Parte Experimental.
"@
$metodoText = @"
A avaliação da disciplina será feita por meio de avaliações escritas individuais (provas) e avaliações de atividades em grupo (relatórios das aulas práticas e/ou trabalhos escritos e/ou apresentações de seminários).
"@
$criterioText = @"
A Média Final (MF) será calculada pela média entre todas as avaliações realizadas durante o semestre, sendo o conjunto das avaliações individuais correspondentes a 75% da composição de MF e o conjunto das avaliações em grupo correspondentes a 25% da composição de MF. Será aprovado o aluno que obtiver MF maior ou igual a cinco e frequência mínima de 70% no semestre.
"@
$normaRecupText = @"
No período de Recuperação haverá horário previamente definido para resolução de dúvidas e será realizada uma avaliação escrita individual (Prova da Recuperação = PR), com conteúdo de todos os tópicos apresentados na disciplina durante o semestre.
A Nota de Recuperação (NR) será dada pela média aritmética entre a Média do Semestre (MF) e a Prova da Recuperação (PR), sendo considerado aprovado o aluno que obtiver NR maior ou igual a cinco.
"@
$bibliografiaText = @"
1) Skoog, D.A.; Holler, F.J. ; Nieman, T.A. Princípios de análise instrumental. 5. ed. Porto Alegre: Bookman,  2002.
2) MENDHAM,J.; DENNEY, R.C.; BARNES, J.D. ; Thomas, M. Vogel: análise química quantitativa. 6. ed. Rio de Janeiro: Livros Técnicos e Científicos, 2002.
3) OHLWEILER, O.A. Fundamentos de análise instrumental. Rio de Janeiro: Livros Técnicos e Científicos, 1981.
4) KRUG, F.J. (org.) Métodos de preparo de amostras: fundamentos sobre métodos de preparo de amostras orgânicas e inorgânicas para análise elementar. 1. ed. Piracicaba: Edição do autor, 2008. 
5) COLLINS, C.H.; BRAGA, G.L.; BONATO, P.S. (Org.) Fundamentos de cromatografia. 1. ed. Campinas: Editora da UNICAMP, 2006.
Bibliografia complementar
1) CHRISTIAN, G.D. Analytical chemistry. 4. ed. Nova York: John Wiley & Sons, 1986.
2) DYER, J.R. Aplicação da espectroscopia de absorção aos compostos orgânicos. 1. Reimpressão. São Paulo: Edgard Blucher, 1977.
3) SILVERSTEIN, R.M.; WEBSTER, F.X.; KIEMLE, D.J. Identificação espectrométrica de compostos orgânicos. 7. ed. Rio de Janeiro: Livros Técnicos e Científicos, 2007.
4) WILLARD, H.H.; MERRITE, L.; DEAB, J. Instrumentação analítica. Lisboa: Fundação Calouste Gulbekian,  1989.
"@

# Step 1: insert a new blank row at position 13, shifting rows 13-23 down to 14-24
$ws.Rows("13:13").Insert()

# Step 2: the inserted row copied formatting from the row above into A13; clear it so
# column A has no cell on row 13 at all (matches target layout).
$ws.Range("A13").Style = "Normal"
$ws.Range("A13").Value = $null

# Step 3: give B13/C13 the same body formatting used elsewhere in columns B/C
# (copy format only, so we reuse the existing style indices instead of creating new ones).
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Step 4: fill in the new "Docentes responsaveis" value on the freshly inserted row.
$ws.Range("B13").Value = $docentesValue
$ws.Range("C13").Value = $docentesValue

# Step 5: "Objetivos:" (row 10) gets the full objectives text instead of the misplaced
# professor name.
$ws.Range("B10").Value = $objetivosText
$ws.Range("C10").Value = $objetivosText

# Step 6: "Programa resumido:" (now row 14 after the insert) gets its real short-syllabus text.
$ws.Range("B14").Value = $programaResumidoText
$ws.Range("C14").Value = $programaResumidoText

# Step 7: "Programa:" (now row 16) gets the full program text.
$ws.Range("B16").Value = $programaText
$ws.Range("C16").Value = $programaText

# Step 8: "Metodo:" (now row 19) gets the evaluation-method text.
$ws.Range("B19").Value = $metodoText
$ws.Range("C19").Value = $metodoText

# Step 9: "Criterio:" (now row 20) gets the grading-criteria text.
$ws.Range("B20").Value = $criterioText
$ws.Range("C20").Value = $criterioText

# Step 10: "Norma de recuperacao:" (now row 21) gets the makeup-exam rule text.
$ws.Range("B21").Value = $normaRecupText
$ws.Range("C21").Value = $normaRecupText

# Step 11: "Bibliografia:" (now row 22) gets the bibliography text.
$ws.Range("B22").Value = $bibliografiaText
$ws.Range("C22").Value = $bibliografiaText
